$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B9").Value = "MPU I2C: "
$ws.Range("C9").Value = "Arduino Guide for MPU-6050 Accelerometer and Gyroscope | Random Nerd Tutorials"
